# Generate Report for Handoff
# Adds a new file entry (37559560-d596-49a6-a40e-260191fa45f7) to the
# localization-status workbook, mirroring the existing
# 217fdef2-d980-41b9-b46a-559cc98513e5 entry, as row 3 on all three sheets.

$wb = $excel.ActiveWorkbook

$newBase = "37559560-d596-49a6-a40e-260191fa45f7"
$newHash = "dcd8221c5d72e27221ccf101d9b2507406536267"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/f92cd8b9dab26ca8166f4125d409ddc687e7dc6a/e2e/$newBase.md"
$zhcnXlf = "$newBase.$newHash.zh-cn.xlf"
$dedeXlf = "$newBase.$newHash.de-de.xlf"
$zhcnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08137d62e8cb7ff71ec9cb4bbc70132f73234939/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhcnXlf"
$dedeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea63a3125bf35593f3bff3b08c559e1f8ebdbdc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$dedeXlf"

# ---------------------------------------------------------------
# Sheet "Overview" -> row 3
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = "$newBase.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-35-21 04:35:26"

$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdUrl, "", "", "$newBase.md")

# ---------------------------------------------------------------
# Sheet "zh-cn" -> row 3
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = "$newBase.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = $zhcnXlf
$ws2.Range("E3").Value = "2016-03-21 04:35:23"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdUrl, "", "", "$newBase.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), $mdUrl, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhcnUrl, "", "", $zhcnXlf)

# ---------------------------------------------------------------
# Sheet "de-de" -> row 3
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = "$newBase.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = $dedeXlf
$ws3.Range("E3").Value = "2016-03-21 04:35:26"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdUrl, "", "", "$newBase.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), $mdUrl, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), $dedeUrl, "", "", $dedeXlf)
